$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: financial-period headers (shift one period newer; 1396/12 drops off, 1401/12 joins) ---
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# --- Row 9: publish dates (same shift) ---
$ws.Range("D9").Value = "1399-03-21 (10)"
$ws.Range("E9").Value = "1400-02-31 (10)"
$ws.Range("F9").Value = "1401-02-31 (11)"
$ws.Range("G9").Value = "1402-02-25 (10)"
$ws.Range("H9").Value = "1402-02-25 (2)"

# --- Row 11: Sales ---
$ws.Range("D11").Value = 27799
$ws.Range("E11").Value = 40364
$ws.Range("F11").Value = 35753
$ws.Range("G11").Value = 48603
$ws.Range("H11").Value = 56612

# --- Row 12: Cost of goods sold ---
$ws.Range("D12").Value = -16505
$ws.Range("E12").Value = -24050
$ws.Range("F12").Value = -16762
$ws.Range("G12").Value = -23531
$ws.Range("H12").Value = -26760

# --- Row 13: Gross profit (loss) ---
$ws.Range("D13").Value = 11293
$ws.Range("E13").Value = 16314
$ws.Range("F13").Value = 18991
$ws.Range("G13").Value = 25072
$ws.Range("H13").Value = 29852

# --- Row 14: General, administrative and organizational expenses ---
$ws.Range("D14").Value = -2699
$ws.Range("E14").Value = -7554
$ws.Range("F14").Value = -4558
$ws.Range("G14").Value = -6744
$ws.Range("H14").Value = -7422

# --- Row 15: Impairment expense - unchanged ("-" every period) ---

# --- Row 16: Net other operating income (expenses) ---
$ws.Range("D16").Value = 409
$ws.Range("E16").Value = 572
$ws.Range("F16").Value = 1281
$ws.Range("G16").Value = 629
$ws.Range("H16").Value = 4920

# --- Row 17: Operating profit (loss) ---
$ws.Range("D17").Value = 9003
$ws.Range("E17").Value = 9332
$ws.Range("F17").Value = 15714
$ws.Range("G17").Value = 18957
$ws.Range("H17").Value = 27349

# --- Row 18: Financial expenses ---
$ws.Range("D18").Value = -2265
$ws.Range("E18").Value = -1050
$ws.Range("F18").Value = -233
$ws.Range("G18").Value = -17
$ws.Range("H18").Value = -52

# --- Row 19: Net other non-operating income and expenses ---
$ws.Range("D19").Value = 458
$ws.Range("E19").Value = 193
$ws.Range("F19").Value = 1580
$ws.Range("G19").Value = 1304
$ws.Range("H19").Value = 2391

# --- Row 20: Net profit (loss) from continuing operations before tax ---
$ws.Range("D20").Value = 7197
$ws.Range("E20").Value = 8475
$ws.Range("F20").Value = 17061
$ws.Range("G20").Value = 20243
$ws.Range("H20").Value = 29688

# --- Row 21: Tax ---
$ws.Range("D21").Value = -330
$ws.Range("E21").Value = -457
$ws.Range("F21").Value = -770
$ws.Range("G21").Value = -422
$ws.Range("H21").Value = -1165

# --- Row 22: Net profit (loss) from continuing operations ---
$ws.Range("D22").Value = 6867
$ws.Range("E22").Value = 8017
$ws.Range("F22").Value = 16290
$ws.Range("G22").Value = 19821
$ws.Range("H22").Value = 28523

# --- Row 23: Discontinued operations profit (loss) - unchanged ("-" every period) ---

# --- Row 24: Net profit (loss) ---
$ws.Range("D24").Value = 6867
$ws.Range("E24").Value = 8017
$ws.Range("F24").Value = 16290
$ws.Range("G24").Value = 19821
$ws.Range("H24").Value = 28523

# --- Row 25: EPS after tax - unchanged (0 every period) ---

# --- Row 26: Capital ---
$ws.Range("D26").Value = 6426
$ws.Range("E26").Value = 10913
$ws.Range("F26").Value = 6192
$ws.Range("G26").Value = 5306
$ws.Range("H26").Value = 3967

# --- Row 27: EPS based on latest capital - unchanged (0 every period) ---
